$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mainImageUrl value in C2 to point to the new logo image
$ws.Range("C2").Value = "https://rmoosa2014.github.io/Resume/LS LOGO.png"

# Update the active selection to match the saved view state (J25)
$ws.Range("J25").Select()
